$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.432.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.62%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.642.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.62%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'196.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +6.85%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'582.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.31%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'3.637.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.65%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.621"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.41%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.14%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +1.53%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +6.71%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'56.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.58%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.0000292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +15.40%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'10.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.49%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'4.227.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.95%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.645.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.92%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +0.57%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'12.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.06%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'68.404.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.80%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +1.56%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +2.30%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'402.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.62%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +26.12%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -1.40%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'86.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.50%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +3.29%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'12.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.25%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +7.14%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'6.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.94%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +19.83%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  +2.43%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'31.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.20%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'706.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +18.07%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'12.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.91%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +5.65%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'64.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.68%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'42.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.51%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.425"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +13.09%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.09%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +6.36%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +19.19%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +4.09%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'3.217.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +19.53%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +13.20%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  +35.51%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  -0.03%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.0422"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.59%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'Stellar"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.133"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.15%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'THORChain"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'8.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +7.83%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'3.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.84%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'Monero"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'142.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.26%  "
$ws.Range("E51").Style = "Normal"

